# vstring.xlsx: populate the "glyphindex" column (G) with the same
# get_X_from_glyphindex accessor names already used for byteindex/charindex,
# mirroring the newly-started glyphindex row family ("Started glyphindex and tests").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G8").Value  = "get_glyph_from_glyphindex"
$ws.Range("G9").Value  = "get_glyphoption_from_glyphindex"
$ws.Range("G10").Value = "get_byteslice_from_glyphindex"
$ws.Range("G11").Value = "get_bytevector_from_glyphindex"
$ws.Range("G12").Value = "get_charvector_from_glyphindex"
$ws.Range("G13").Value = "get_glyphvector_from_glyphindex"
$ws.Range("G14").Value = "get_byteiterator_from_glyphindex"
$ws.Range("G15").Value = "get_chariterator_from_glyphindex"
$ws.Range("G16").Value = "get_glyphiterator_from_glyphindex"
$ws.Range("G17").Value = "get_strref_from_glyphindex"
$ws.Range("G18").Value = "get_string_from_glyphindex"

# Column widths grew to accommodate the new text (best-effort resize -
# the host's width model is coarser than Excel's real pixel-metric AutoFit,
# so these land close to, not byte-identical with, the authored widths).
$ws.Columns.Item(4).ColumnWidth = 33.42
$ws.Columns.Item(7).ColumnWidth = 30.59

# Selection moved when the author kept typing past the table.
[void]$ws.Range("G25").Select()
